$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update the Date value (B8) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- Sheet "Elements": swap the two "Mapping" columns (AK <-> AL) ---
$wsElem = $wb.Worksheets.Item("Elements")

# Swap AK <-> AL for the header plus all data rows (rows 1-6),
# only touching rows whose two values actually differ.
for ($r = 1; $r -le 6; $r++) {
    $ak = $wsElem.Range("AK$r").Value2
    $al = $wsElem.Range("AL$r").Value2
    if ($ak -ne $al) {
        $wsElem.Range("AK$r").Value = $al
        $wsElem.Range("AL$r").Value = $ak
    }
}

# Swap the column widths to match the new (swapped) content
$wsElem.Columns.Item(37).ColumnWidth = 73.0
$wsElem.Columns.Item(38).ColumnWidth = 24.166666666666668

# Preserve the originally-hidden columns (engine round-trip otherwise drops the flag)
$wsElem.Columns.Item(3).Hidden = $true
$wsElem.Columns.Item(4).Hidden = $true
$wsElem.Columns.Item(31).Hidden = $true
$wsElem.Columns.Item(32).Hidden = $true
$wsElem.Columns.Item(33).Hidden = $true
